# Rectification de code de config
# - Fixes product codes/names that were mismatched against their quantities/amounts
#   (several rows had the wrong ARTICLE/DESIGNATION paired with QTE/TTC values).
# - Renames the QTE/TTC headers to their full names QUANTITE / MONTANT TTC.
# - Widens the QUANTITE/MONTANT TTC columns to fit the new header text.
# - Recomputes the TOTAUX row to match the corrected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (C: QUANTITE, D: MONTANT TTC) ---
# ColumnWidth is expressed in "characters" and gets pixel-quantized on write,
# so we target the input that rounds to the desired stored width (12 and ~15.6).
$ws.Columns.Item(3).ColumnWidth = 11.166666666666668
$ws.Columns.Item(4).ColumnWidth = 14.833333333333332

# Row 1
$ws.Cells.Item(1, 3).Value = "QUANTITE"
$ws.Cells.Item(1, 4).Value = "MONTANT TTC"
# Row 3
$ws.Cells.Item(3, 1).Value = "BAH022"
$ws.Cells.Item(3, 2).Value = "BEEF FILET/ FILET DE ZEBU"
$ws.Cells.Item(3, 3).Value = 13.556
$ws.Cells.Item(3, 4).Value = 776216.56
# Row 4
$ws.Cells.Item(4, 1).Value = "BURG0001"
$ws.Cells.Item(4, 2).Value = "BURGER 100% VIANDE / PACK OF 2"
$ws.Cells.Item(4, 3).Value = 13.76
$ws.Cells.Item(4, 4).Value = 750415.35
# Row 6
$ws.Cells.Item(6, 1).Value = "BAF008"
$ws.Cells.Item(6, 2).Value = "RIBEYE / ENTRECOTE"
$ws.Cells.Item(6, 3).Value = 12.848
$ws.Cells.Item(6, 4).Value = 713192.48
# Row 7
$ws.Cells.Item(7, 1).Value = "BLPCG001"
$ws.Cells.Item(7, 2).Value = "BLANC DE POULET CONGELE"
$ws.Cells.Item(7, 3).Value = 23.282
$ws.Cells.Item(7, 4).Value = 675178.0
# Row 11
$ws.Cells.Item(11, 1).Value = "MERL0001"
$ws.Cells.Item(11, 2).Value = " LAMB MERGUEZ / MERGUEZ D'AGNEAU "
$ws.Cells.Item(11, 3).Value = 6.302
$ws.Cells.Item(11, 4).Value = 383665.76
# Row 12
$ws.Cells.Item(12, 1).Value = "BEESK001"
$ws.Cells.Item(12, 2).Value = "BEEF SKEWERS - 100% FILET - PREMIUM"
$ws.Cells.Item(12, 3).Value = 5.84
$ws.Cells.Item(12, 4).Value = 365058.4
# Row 20
$ws.Cells.Item(20, 1).Value = "CUPCG001"
$ws.Cells.Item(20, 2).Value = "CUISSE ENTIÈRE DE POULET CONGELE"
$ws.Cells.Item(20, 3).Value = 6.61
$ws.Cells.Item(20, 4).Value = 165250.0
# Row 21
$ws.Cells.Item(21, 1).Value = "GOAR0002"
$ws.Cells.Item(21, 2).Value = "GOAT RIB / COTE DE CHEVRE"
$ws.Cells.Item(21, 3).Value = 2.814
$ws.Cells.Item(21, 4).Value = 154910.7
# Row 22
$ws.Cells.Item(22, 1).Value = "BOUC0001"
$ws.Cells.Item(22, 2).Value = " BOURGUIGNON CUBE"
$ws.Cells.Item(22, 3).Value = 3.88
$ws.Cells.Item(22, 4).Value = 153337.6
# Row 23
$ws.Cells.Item(23, 1).Value = "PLPCG001"
$ws.Cells.Item(23, 2).Value = "PILON CONGELE"
$ws.Cells.Item(23, 3).Value = 6.016
$ws.Cells.Item(23, 4).Value = 150400.0
# Row 24
$ws.Cells.Item(24, 1).Value = "HCPCG001"
$ws.Cells.Item(24, 2).Value = "HAUT DE CUISSE CONGELE"
$ws.Cells.Item(24, 3).Value = 5.612
$ws.Cells.Item(24, 4).Value = 140300.0
# Row 25
$ws.Cells.Item(25, 1).Value = "ALPCG001"
$ws.Cells.Item(25, 2).Value = "AILES DE POULET CONGELE"
$ws.Cells.Item(25, 3).Value = 5.626
$ws.Cells.Item(25, 4).Value = 137837.0
# Row 28
$ws.Cells.Item(28, 1).Value = "GOAS0001"
$ws.Cells.Item(28, 2).Value = "GOAT STEW / CHEVRE CUBE AVEC OS"
$ws.Cells.Item(28, 3).Value = 1.952
$ws.Cells.Item(28, 4).Value = 96477.6
# Row 29
$ws.Cells.Item(29, 1).Value = "FOICG001"
$ws.Cells.Item(29, 2).Value = "FOIE CONGELE"
$ws.Cells.Item(29, 3).Value = 5.492
$ws.Cells.Item(29, 4).Value = 96110.0
# Row 30
$ws.Cells.Item(30, 1).Value = "RUMT0001"
$ws.Cells.Item(30, 2).Value = "RUMP TAIL / AIGUILLETTE"
$ws.Cells.Item(30, 3).Value = 1.204
$ws.Cells.Item(30, 4).Value = 79328.56
# Row 42
$ws.Cells.Item(42, 1).Value = "MERC0001"
$ws.Cells.Item(42, 2).Value = "MERGUEZ CONGELE"
$ws.Cells.Item(42, 3).Value = 1.292
$ws.Cells.Item(42, 4).Value = 38760.0
# Row 43
$ws.Cells.Item(43, 1).Value = "JMBVOL04"
$ws.Cells.Item(43, 2).Value = "JAMBON DE VOLLAILE MALAGASY CONGELE"
$ws.Cells.Item(43, 3).Value = 1.074
$ws.Cells.Item(43, 4).Value = 35979.0
# Row 44
$ws.Cells.Item(44, 1).Value = "BEER0001"
$ws.Cells.Item(44, 2).Value = "BEEF RIBS / PLAT DE COTE"
$ws.Cells.Item(44, 3).Value = 0.81
$ws.Cells.Item(44, 4).Value = 35696.7
# Row 45
$ws.Cells.Item(45, 1).Value = "PEFR0001"
$ws.Cells.Item(45, 2).Value = "POULET ENTIER FRAIS"
$ws.Cells.Item(45, 3).Value = 1.34
$ws.Cells.Item(45, 4).Value = 33768.0
# Row 46
$ws.Cells.Item(46, 1).Value = "BEES0001"
$ws.Cells.Item(46, 2).Value = "BEEF STRIRFY / EMINCE DE ZEBU"
$ws.Cells.Item(46, 3).Value = 0.708
$ws.Cells.Item(46, 4).Value = 29130.66
# Row 47
$ws.Cells.Item(47, 1).Value = "SACHET_0"
$ws.Cells.Item(47, 2).Value = "PLASTIC BAGS"
$ws.Cells.Item(47, 3).Value = 41.0
$ws.Cells.Item(47, 4).Value = 20500.0
# Row 48
$ws.Cells.Item(48, 1).Value = "BRIA0001"
$ws.Cells.Item(48, 2).Value = "AMERICAN BRISKET"
$ws.Cells.Item(48, 3).Value = 0.322
$ws.Cells.Item(48, 4).Value = 15790.88
# Row 49
$ws.Cells.Item(49, 1).Value = "BAH005"
$ws.Cells.Item(49, 2).Value = "TOPSIDE / GITE DE ZEBU"
$ws.Cells.Item(49, 3).Value = 0.352
$ws.Cells.Item(49, 4).Value = 15597.12
# Row 50
$ws.Cells.Item(50, 3).Value = 292.184
$ws.Cells.Item(50, 4).Value = 10389012.89

Write-Output "Palmares produit depuis: corrections applied."
